$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 18 (shifts existing rows 18-39 down to 19-40),
# preserving the formatting of the row being pushed down.
$ws.Rows.Item(18).Insert()

# Populate the newly inserted row 18 with this week's data.
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 44966
$ws.Cells.Item(18, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = 100112010
$ws.Cells.Item(18, 7).Value = "Achicoria"
$ws.Cells.Item(18, 8).Value = "Sin especificar"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 70
$ws.Cells.Item(18, 11).Value = 7000
$ws.Cells.Item(18, 12).Value = 7000
$ws.Cells.Item(18, 13).Value = 7000
$ws.Cells.Item(18, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(18, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(18, 16).Value = 438
$ws.Cells.Item(18, 17).Value = 16
$ws.Cells.Item(18, 18).Value = "Hortaliza"
